# Controls.xlsx — "Added a new control."
#
# Two new data rows are introduced on Blad1 (the control-mapping table,
# columns A=WindowClassName, B=windowControlID, C=Module, D=Text, E=Comment):
#
#   1) A brand-new row is inserted at row 52 (pushing the former rows
#      52-331 down to 53-332):
#         A52 = "Edit"            (existing shared string)
#         B52 = 20593
#         D52 = "Anteckningar"    (new shared string)
#
#   2) A brand-new row is appended at the end, row 333:
#         A333 = "ComboBox"                                        (existing shared string)
#         B333 = 25019
#         D333 = "Inställningar för egen kopia, utskriftsval"      (new shared string)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) insert the new "Anteckningar" control row at row 52 ---------------
$ws.Rows.Item(52).Insert()

$ws.Range("A52").Value = "Edit"
$ws.Range("B52").Value = 20593
$ws.Range("D52").Value = "Anteckningar"

# --- 2) append the new "Inställningar för egen kopia, utskriftsval" row ---
$ws.Range("A333").Value = "ComboBox"
$ws.Range("B333").Value = 25019
$ws.Range("D333").Value = "Inställningar för egen kopia, utskriftsval"

# --- reflect the new selection / scroll position left by the edit ---------
$ws.Range("D333").Select()
